# Generate Report for Handoff
#
# A new handoff round-trip just completed for the
# "923d649d-0171-4c9d-a1f2-32ac1819f07c.md" source file, so the
# "Latest Handoff Datetime" column (D) on row 4 of each language sheet
# is refreshed with the new handoff timestamp. Nothing else on the
# report changes.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-17 05:51:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-17 05:51:59"
